$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.393.10'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.937.08'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''0.7682'
$ws.Range('E5').Value = '  +8.33%  '
$ws.Range('D6').Value = '''248.10'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '''27.91'
$ws.Range('D9').Value = '''0.3198'
$ws.Range('E9').Value = '  -3.17%  '
$ws.Range('D10').Value = '''0.07097'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('D11').Value = '''0.7820'
$ws.Range('E11').Value = '  -2.85%  '
$ws.Range('D12').Value = '''0.08009'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').Value = '1.936.47'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '''5.373'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').Value = '''95.03'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').Value = '''14.52'
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('D17').Value = '30.391.60'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '''257.60'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '''0.000008007'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').Value = '''5.850'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '2.194.00'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '''6.757'
$ws.Range('E24').Value = '  -3.51%  '
$ws.Range('D25').Value = '''9.601'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('D26').Value = '''164.42'
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('D27').Value = '''19.15'
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('D28').Value = '''0.1331'
$ws.Range('E28').Value = '  +2.85%  '
$ws.Range('D29').Value = '''2.290'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').Value = '''1.364'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').Value = '''1.528'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').Value = '''4.435'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('D33').Value = '''4.150'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').Value = '''0.05191'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('D35').Value = '''1.278'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('D36').Value = '''0.7493'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '''2.780'
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').Value = '''0.01970'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').Value = '''78.29'
$ws.Range('D41').Value = '''6.461'
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').Value = '''0.4517'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = '''1.972'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').Value = '''0.8349'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').Value = '''101.23'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').Value = '''9.790'
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').Value = '''7.515'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''985.49'
$ws.Range('E49').Value = '  +11.35%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''37.45'
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '''0.4158'
$ws.Range('E51').Value = '  -0.51%  '
